$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Förändrad" (C column) date from 2023-09-14 (45183) to 2023-09-15 (45184) for rows 2-81 ---
for ($r = 2; $r -le 81; $r++) {
    $ws.Range("C$r").Value = 45184
}

# --- Rebuild HYPERLINK formulas with an explicit display-text second argument ---

$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/artfynd/A 31572-2023.xlsx, "A 31572-2023"")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/kartor/A 31572-2023.png", "A 31572-2023")'
$ws.Range("U2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/knärot/A 31572-2023.png", "A 31572-2023")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/klagomål/A 31572-2023.docx", "A 31572-2023")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/klagomålsmail/A 31572-2023.docx", "A 31572-2023")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/tillsyn/A 31572-2023.docx", "A 31572-2023")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/tillsynsmail/A 31572-2023.docx", "A 31572-2023")'
$ws.Range("S3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/artfynd/A 32292-2023.xlsx, "A 32292-2023"")'
$ws.Range("T3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/kartor/A 32292-2023.png", "A 32292-2023")'
$ws.Range("V3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/klagomål/A 32292-2023.docx", "A 32292-2023")'
$ws.Range("W3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/klagomålsmail/A 32292-2023.docx", "A 32292-2023")'
$ws.Range("X3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/tillsyn/A 32292-2023.docx", "A 32292-2023")'
$ws.Range("Y3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/tillsynsmail/A 32292-2023.docx", "A 32292-2023")'
$ws.Range("S4").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/artfynd/A 32299-2023.xlsx, "A 32299-2023"")'
$ws.Range("T4").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/kartor/A 32299-2023.png", "A 32299-2023")'
$ws.Range("V4").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/klagomål/A 32299-2023.docx", "A 32299-2023")'
$ws.Range("W4").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/klagomålsmail/A 32299-2023.docx", "A 32299-2023")'
$ws.Range("X4").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/tillsyn/A 32299-2023.docx", "A 32299-2023")'
$ws.Range("Y4").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/tillsynsmail/A 32299-2023.docx", "A 32299-2023")'
$ws.Range("S5").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/artfynd/A 32785-2023.xlsx, "A 32785-2023"")'
$ws.Range("T5").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/kartor/A 32785-2023.png", "A 32785-2023")'
$ws.Range("V5").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/klagomål/A 32785-2023.docx", "A 32785-2023")'
$ws.Range("W5").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/klagomålsmail/A 32785-2023.docx", "A 32785-2023")'
$ws.Range("X5").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/tillsyn/A 32785-2023.docx", "A 32785-2023")'
$ws.Range("Y5").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/tillsynsmail/A 32785-2023.docx", "A 32785-2023")'
$ws.Range("U33").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/knärot/A 31569-2023.png", "A 31569-2023")'
$ws.Range("V33").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/klagomål/A 31569-2023.docx", "A 31569-2023")'
$ws.Range("W33").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/klagomålsmail/A 31569-2023.docx", "A 31569-2023")'
$ws.Range("X33").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/tillsyn/A 31569-2023.docx", "A 31569-2023")'
$ws.Range("Y33").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HEBY/tillsynsmail/A 31569-2023.docx", "A 31569-2023")'
